$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").NumberFormat = "@"

$ws.Range("E2").Value = "2026-02-20 10:10:58"
$ws.Range("H2").Value = "91%"
$ws.Range("J2").Value = "1020.8 hPa"
$ws.Range("K2").Value = "3.0 MJ/m2"
$ws.Range("M2").Value = "13.5 °C 9:45 TU"
$ws.Range("O2").Value = "3.3 °C"
